$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 120 (pushes the old rows 120-125 down to 122-127).
$ws.Rows("120:121").Insert()

# --- New row 120: Damasco / Modesto / Primera, week of 2023-01-20 ---
$ws.Cells.Item(120, 1).Value = 8
$ws.Cells.Item(120, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(120, 3).Value = "Coquimbo"
$ws.Cells.Item(120, 4).Value = 44946
$ws.Cells.Item(120, 5).Value = 4
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100103
$ws.Cells.Item(120, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(120, 9).Value = 100103003
$ws.Cells.Item(120, 10).Value = "Damasco"
$ws.Cells.Item(120, 11).Value = "Modesto"
$ws.Cells.Item(120, 12).Value = "Primera"
$ws.Cells.Item(120, 13).Value = 160
$ws.Cells.Item(120, 14).Value = 20000
$ws.Cells.Item(120, 15).Value = 21000
$ws.Cells.Item(120, 16).Value = 20500
$ws.Cells.Item(120, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(120, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(120, 19).Value = 1281
$ws.Cells.Item(120, 20).Value = 16

# --- New row 121: Damasco / Modesto / Segunda, week of 2023-01-20 ---
$ws.Cells.Item(121, 1).Value = 8
$ws.Cells.Item(121, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(121, 3).Value = "Coquimbo"
$ws.Cells.Item(121, 4).Value = 44946
$ws.Cells.Item(121, 5).Value = 4
$ws.Cells.Item(121, 6).Value = "Fruta"
$ws.Cells.Item(121, 7).Value = 100103
$ws.Cells.Item(121, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(121, 9).Value = 100103003
$ws.Cells.Item(121, 10).Value = "Damasco"
$ws.Cells.Item(121, 11).Value = "Modesto"
$ws.Cells.Item(121, 12).Value = "Segunda"
$ws.Cells.Item(121, 13).Value = 140
$ws.Cells.Item(121, 14).Value = 16000
$ws.Cells.Item(121, 15).Value = 17000
$ws.Cells.Item(121, 16).Value = 16500
$ws.Cells.Item(121, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(121, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(121, 19).Value = 1031
$ws.Cells.Item(121, 20).Value = 16
